$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.406.76"
$ws.Range("E2").Value = "  +3.19%  "

$ws.Range("D3").Value = "1.871.92"
$ws.Range("E3").Value = "  +1.35%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.53%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.41"
$ws.Range("E5").Value = "  +1.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.58%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4689"
$ws.Range("E7").Value = "  +1.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3961"
$ws.Range("E8").Value = "  +2.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.84"
$ws.Range("E9").Value = "  +3.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08030"
$ws.Range("E10").Value = "  +1.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.003"
$ws.Range("E11").Value = "  +1.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.90"
$ws.Range("E12").Value = "  +2.80%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.992"
$ws.Range("E13").Value = "  +1.19%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.860.58"
$ws.Range("E14").Value = "  +0.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.245"
$ws.Range("E15").Value = "  +2.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.06"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001040"
$ws.Range("E18").Value = "  +0.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06612"
$ws.Range("E19").Value = "  -0.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.58"
$ws.Range("E20").Value = "  +3.17%  "

$ws.Range("E21").Value = "  -0.55%  "

$ws.Range("D22").Value = "28.410.40"
$ws.Range("E22").Value = "  +3.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.447"
$ws.Range("E23").Value = "  +1.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.05"
$ws.Range("E24").Value = "  +1.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.272"
$ws.Range("E25").Value = "  -1.23%  "

$ws.Range("D26").Value = "2.093.81"
$ws.Range("E26").Value = "  +1.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.29"
$ws.Range("E27").Value = "  +1.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.75"
$ws.Range("E28").Value = "  +1.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.138"
$ws.Range("E29").Value = "  +2.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.514"
$ws.Range("E30").Value = "  +2.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.07"
$ws.Range("E31").Value = "  +0.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9711"
$ws.Range("E32").Value = "  +1.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09489"
$ws.Range("E33").Value = "  +1.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.570"
$ws.Range("E34").Value = "  -0.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.380"
$ws.Range("E35").Value = "  +3.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.352"
$ws.Range("E36").Value = "  +1.64%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06110"
$ws.Range("E37").Value = "  +2.70%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02257"
$ws.Range("E38").Value = "  +2.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.377"
$ws.Range("E39").Value = "  +3.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.184"
$ws.Range("E40").Value = "  +1.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5952"
$ws.Range("E41").Value = "  +1.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  -0.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1874"
$ws.Range("E43").Value = "  +1.23%  "

$ws.Range("E44").Value = "  +2.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.295"
$ws.Range("E45").Value = "  +3.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5576"
$ws.Range("E46").Value = "  +0.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.20"
$ws.Range("E47").Value = "  +1.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.960"
$ws.Range("E48").Value = "  +4.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06864"
$ws.Range("E49").Value = "  +2.80%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.047"
$ws.Range("E50").Value = "  +14.17%  "

$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.33"
$ws.Range("E51").Value = "  +0.67%  "
